$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Range("D21").Value = "Intake,Supervisor Approval"
$ws.Range("D22").Value = "Fulfill,Executive Approval"
$ws.Range("D23").Value = "Supervisor Approval,Release"
$ws.Range("D24").Value = "Executive Approval"
$ws.Range("E24").Value = ""
